# Apply the workbook update:
#  1. Column C ("Förändrad") for every data row (2..224): 45184 -> 45186
#  2. For rows 2..19, every HYPERLINK(...) formula in columns S/T/U/V/W/X/Y
#     gets a second argument added: the friendly-name text (column A's
#     "Beteckning" value for that row), e.g.
#       HYPERLINK("...xlsx")  ->  HYPERLINK("...xlsx", "A 49304-2019")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow  = 224
$oldDate = 45184
$newDate = 45186

$hyperlinkCols = @("S", "T", "U", "V", "W", "X", "Y")

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {

    # --- 1) Bump the "Förändrad" date in column C, if present ---
    $cCell = $ws.Cells.Item($r, 3)
    $cVal = $cCell.Value2
    if ($cVal -eq $oldDate) {
        $cCell.Value = $newDate
    }

    # --- 2) Add the display-text argument to HYPERLINK formulas ---
    $label = $ws.Cells.Item($r, 1).Value2
    if ($label) {
        foreach ($col in $hyperlinkCols) {
            $cell = $ws.Range("$col$r")
            $f = $cell.Formula
            if ($f -and $f.Length -gt 0 -and $f -like '*HYPERLINK(*' -and $f -notlike '*,*') {
                $newF = $f.Substring(0, $f.Length - 1) + ', "' + $label + '")'
                $cell.Formula = $newF
            }
        }
    }
}
